# #3456 updated PM Property ID
# Update the "Portfolio Manager Building ID" values (column B, rows 2-10)
# on the active sheet to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newIds = @(22178843, 22178844, 22178845, 22178846, 22178847, 22178848, 22178849, 22178850, 22178851)

$row = 2
foreach ($id in $newIds) {
    $ws.Cells.Item($row, 2).Value = $id
    $row = $row + 1
}

# Match the selection left behind in the saved workbook (B2:B10, active cell B2).
$ws.Range("B2:B10").Select()
